$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 576.6667
$ws.Range("I32").Value = 375
$ws.Range("K32").Value = 375
$ws.Range("M32").Value = -49
$ws.Range("H40").Value = 1958.25
$ws.Range("I40").Value = 1750
$ws.Range("J40").Value = 2166.5
$ws.Range("K40").Value = 1750
$ws.Range("L40").Value = 2166.5
$ws.Range("M40").Value = -1575
$ws.Range("N40").Value = -2516.5
$ws.Range("H51").Value = 8173.3335
$ws.Range("J51").Value = 11260
$ws.Range("L51").Value = 11260
$ws.Range("N51").Value = -12228
$ws.Range("H111").Value = 5048.6
$ws.Range("I111").Value = 1797.2
$ws.Range("K111").Value = 5391.6
$ws.Range("M111").Value = -2324.6
$ws.Range("H112").Value = 2778789.2
$ws.Range("I112").Value = 350
$ws.Range("J112").Value = 2925023
$ws.Range("K112").Value = 1050
$ws.Range("L112").Value = 8775069
$ws.Range("M112").Value = 58
$ws.Range("N112").Value = -8777285
$ws.Range("H113").Value = 20411202
$ws.Range("I113").Value = 31251950
$ws.Range("J113").Value = 5088.1763
$ws.Range("K113").Value = 31251950
$ws.Range("L113").Value = 5088.1763
$ws.Range("M113").Value = -31248696
$ws.Range("N113").Value = -11596.1763
$ws.Range("H129").Value = 847.49054
$ws.Range("I129").Value = 795
$ws.Range("J129").Value = 849.549
$ws.Range("K129").Value = 2385
$ws.Range("L129").Value = 2548.647
$ws.Range("M129").Value = 2615
$ws.Range("N129").Value = -12548.647
$ws.Range("H137").Value = 32320.727
$ws.Range("I137").Value = 2063.65
$ws.Range("K137").Value = 6190.950000000001
$ws.Range("M137").Value = -3640.950000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2015.7391
$ws.Range("I2").Value = 1287.5625
$ws.Range("J2").Value = 3680.1428
$ws.Range("K2").Value = 1287.5625
$ws.Range("L2").Value = 3680.1428
$ws.Range("M2").Value = -1174.5625
$ws.Range("N2").Value = -3906.1428
$ws.Range("H32").Value = 26168.285
$ws.Range("I32").Value = 27154.574
$ws.Range("K32").Value = 27154.574
$ws.Range("M32").Value = -26867.574
$ws.Range("H45").Value = 3162.7856
$ws.Range("I45").Value = 3048.4119
$ws.Range("J45").Value = 3240.56
$ws.Range("K45").Value = 3048.4119
$ws.Range("L45").Value = 3240.56
$ws.Range("M45").Value = -2671.4119
$ws.Range("N45").Value = -3994.56
$ws.Range("H61").Value = 3529.9524
$ws.Range("I61").Value = 2891.8333
$ws.Range("J61").Value = 4380.778
$ws.Range("K61").Value = 2891.8333
$ws.Range("L61").Value = 4380.778
$ws.Range("M61").Value = -2679.8333
$ws.Range("N61").Value = -4804.778
$ws.Range("H97").Value = 1538.4445
$ws.Range("I97").Value = 1448.2667
$ws.Range("K97").Value = 1448.2667
$ws.Range("M97").Value = -952.2666999999999
$ws.Range("H102").Value = 2188.1304
$ws.Range("I102").Value = 821.75
$ws.Range("K102").Value = 821.75
$ws.Range("M102").Value = 800.25
$ws.Range("H116").Value = 2015.7391
$ws.Range("I116").Value = 1287.5625
$ws.Range("J116").Value = 3680.1428
$ws.Range("K116").Value = 1287.5625
$ws.Range("L116").Value = 3680.1428
$ws.Range("M116").Value = 1006.4375
$ws.Range("N116").Value = -8268.1428
$ws.Range("H136").Value = 3529.9524
$ws.Range("I136").Value = 2891.8333
$ws.Range("J136").Value = 4380.778
$ws.Range("K136").Value = 8675.499899999999
$ws.Range("L136").Value = 13142.334
$ws.Range("M136").Value = -6125.499899999999
$ws.Range("N136").Value = -18242.334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2015.7391
$ws.Range("I3").Value = 1287.5625
$ws.Range("J3").Value = 3680.1428
$ws.Range("K3").Value = 1287.5625
$ws.Range("L3").Value = 3680.1428
$ws.Range("M3").Value = -1173.5625
$ws.Range("N3").Value = -3908.1428
$ws.Range("H22").Value = 188.11111
$ws.Range("I22").Value = 188.11111
$ws.Range("K22").Value = 188.11111
$ws.Range("M22").Value = -15.11111
$ws.Range("H94").Value = 3469.75
$ws.Range("I94").Value = 1638.579
$ws.Range("K94").Value = 1638.579
$ws.Range("M94").Value = -1187.579
$ws.Range("H99").Value = 2188.9
$ws.Range("I99").Value = 1878.7778
$ws.Range("K99").Value = 1878.7778
$ws.Range("M99").Value = -380.7778000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9768.436
$ws.Range("I31").Value = 11066.25
$ws.Range("K31").Value = 11066.25
$ws.Range("M31").Value = -10771.25
$ws.Range("H34").Value = 9768.436
$ws.Range("I34").Value = 11066.25
$ws.Range("K34").Value = 11066.25
$ws.Range("M34").Value = -10864.25
$ws.Range("H62").Value = 5762.3
$ws.Range("I62").Value = 4121
$ws.Range("K62").Value = 4121
$ws.Range("M62").Value = -3497
$ws.Range("H65").Value = 5762.3
$ws.Range("I65").Value = 4121
$ws.Range("K65").Value = 20605
$ws.Range("M65").Value = -17485
$ws.Range("H86").Value = 15706.637
$ws.Range("J86").Value = 22838
$ws.Range("L86").Value = 22838
$ws.Range("N86").Value = -25084
$ws.Range("H89").Value = 15706.637
$ws.Range("J89").Value = 22838
$ws.Range("L89").Value = 114190
$ws.Range("N89").Value = -125422
$ws.Range("H107").Value = 1606.25
$ws.Range("I107").Value = 1244.9
$ws.Range("J107").Value = 3413
$ws.Range("K107").Value = 1244.9
$ws.Range("L107").Value = 3413
$ws.Range("M107").Value = 675.0999999999999
$ws.Range("N107").Value = -7253
$ws.Range("H132").Value = 24098.916
$ws.Range("I132").Value = 31749.941
$ws.Range("K132").Value = 95249.823
$ws.Range("M132").Value = -92719.823

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 738.8
$ws.Range("J131").Value = 754.2708
$ws.Range("L131").Value = 2262.8124
$ws.Range("N131").Value = -12342.8124
$ws.Range("H132").Value = 980
$ws.Range("I132").Value = 980
$ws.Range("K132").Value = 8820
$ws.Range("M132").Value = -6290

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1816.2858
$ws.Range("I22").Value = 2180
$ws.Range("J22").Value = 1614.2222
$ws.Range("K22").Value = 2180
$ws.Range("L22").Value = 1614.2222
$ws.Range("M22").Value = -1885
$ws.Range("N22").Value = -2204.2222
$ws.Range("H27").Value = 1816.2858
$ws.Range("I27").Value = 2180
$ws.Range("J27").Value = 1614.2222
$ws.Range("K27").Value = 2180
$ws.Range("L27").Value = 1614.2222
$ws.Range("M27").Value = -2073
$ws.Range("N27").Value = -1828.2222
$ws.Range("H46").Value = 2523.3333
$ws.Range("J46").Value = 2733.6667
$ws.Range("L46").Value = 2733.6667
$ws.Range("N46").Value = -3109.6667
$ws.Range("H55").Value = 176.52942
$ws.Range("I55").Value = 136.6
$ws.Range("J55").Value = 193.16667
$ws.Range("K55").Value = 136.6
$ws.Range("L55").Value = 193.16667
$ws.Range("M55").Value = 36.40000000000001
$ws.Range("N55").Value = -539.1666700000001
$ws.Range("H61").Value = 4893.75
$ws.Range("I61").Value = 2003.6364
$ws.Range("K61").Value = 2003.6364
$ws.Range("M61").Value = -1801.6364
$ws.Range("H68").Value = 3905.3914
$ws.Range("I68").Value = 1875.091
$ws.Range("J68").Value = 5766.5
$ws.Range("K68").Value = 1875.091
$ws.Range("L68").Value = 5766.5
$ws.Range("M68").Value = -1126.091
$ws.Range("N68").Value = -7264.5
$ws.Range("H71").Value = 3905.3914
$ws.Range("I71").Value = 1875.091
$ws.Range("J71").Value = 5766.5
$ws.Range("K71").Value = 9375.455
$ws.Range("L71").Value = 28832.5
$ws.Range("M71").Value = -5631.455
$ws.Range("N71").Value = -36320.5
$ws.Range("H93").Value = 1777.6666
$ws.Range("I93").Value = 1754.2632
$ws.Range("K93").Value = 1754.2632
$ws.Range("M93").Value = -506.2632000000001
$ws.Range("H113").Value = 4893.75
$ws.Range("I113").Value = 2003.6364
$ws.Range("K113").Value = 2003.6364
$ws.Range("M113").Value = 166.3635999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3379764
$ws.Range("I113").Value = 1780.6666
$ws.Range("K113").Value = 5341.9998
$ws.Range("M113").Value = -3171.9998
